# Added Google Charts Dashboard for Sony Forum
# Rework FinalData.xlsx Sheet1 from the Samsung scraping sample (5 cols x 4 rows)
# into the Sony scraping sample (8 cols x 2 rows).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Drop the extra Samsung M20 rows (rows 3 & 4) - only one data row remains.
# ---------------------------------------------------------------------------
$ws.Rows.Item(4).Delete()
$ws.Rows.Item(3).Delete()

# ---------------------------------------------------------------------------
# 2) Header row (row 1): rename existing headers, add 3 new header cells.
#    Column A (Product) and D (Category) keep their original text.
# ---------------------------------------------------------------------------
$ws.Range("B1").Value = "First User Name"
$ws.Range("C1").Value = "Second User Name"
$ws.Range("E1").Value = "Thread Name"

# New header cells F1:H1 - copy the header style (bold/border/center) from A1
# first, then set their text.
$ws.Range("A1").Copy()
$ws.Range("F1:H1").PasteSpecial(-4122)
$ws.Range("F1").Value = "Links "
$ws.Range("G1").Value = "Date"
$ws.Range("H1").Value = "Issue Detail"

# ---------------------------------------------------------------------------
# 3) Data row (row 2): replace the Samsung Galaxy M30 record with the Sony
#    Xperia 1 "Fast Charger UCH32C" thread, and populate the 3 new columns.
# ---------------------------------------------------------------------------
$ws.Range("A2").Value = "Xperia 1"

$jonas = @"

Jonas

"@
$ws.Range("B2").Value = $jonas

$ws.Range("C2").Value = " Sony Xperia Support"
$ws.Range("D2").Value = "Touch,"
$ws.Range("E2").Value = "Fast Charger UCH32C"
$ws.Range("F2").Value = "https://talk.sonymobile.com/t5/Xperia-1/Fast-Charger-UCH32C/td-p/1370427"

# "2019-04-11" must stay a plain text value (not get reinterpreted as a
# date serial number). Force text with a leading apostrophe, then restore
# plain (non quote-prefixed) formatting by pasting the format from a
# neighbouring plain-text cell.
$ws.Range("G2").Value = "'2019-04-11"
$ws.Range("A2").Copy()
$ws.Range("G2").PasteSpecial(-4122)

$issueDetail = @"

Hi @pressefr!Do I understand correctly if you want to know if this charger will be included in box content for the Xperia 1 in the US? If that is correct, this is very tricky for me to answer here at our Global user based forum as the box content is different, not only between markets, but can also be different depending on the retailer from where the device is bought.The UCH32C charger is not listed on the Sony Mobile website for the US, but you may want to get in touch with the Local support team there and see if they have any additional information, not only regarding the availability, but also the box content for the Xperia 1, in their market.Feel free to let me know if you are able to get some information about this from them, as I'm sure there are other users from the US that are also interested in knowing this! 
Official Sony Xperia Support StaffIf you're new to our forums, make sure that you've read our Discussion guidelines.To get in touch with your local support team, please visit our contact page.

"@
$ws.Range("H2").Value = $issueDetail

# Undo the auto row-height ("wrap"/auto-fit) bump that typing a multi-line
# value into row 2 triggers, so row 2 keeps its default height.
$ws.Rows.Item(2).AutoFit()
